# LumFunc.xlsx update:
#  - Added direct SAGA Quenching Plots (wpd_datasets (3) text-import connection/query table
#    metadata cannot be represented through this COM surface, so the closest achievable
#    equivalent -- the resulting worksheet data, shared strings and the scoped defined
#    name Excel creates for an imported text/csv range -- is reproduced below).
#  - Updated LFs with new MW analogs (NGC4258 / NGC4631 X/Y series in columns M:P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New luminosity-function columns for the two new Milky Way analogs ---

# Headers (row 1) -- only the first column of each X/Y pair carries a label,
# matching the existing M31/M101/M94/CenA/M81 layout in columns A,C,E,G,I,K.
$ws.Range("M1").Value = "NGC4258"
$ws.Range("O1").Value = "NGC4631"

# Column sub-headers (row 2)
$ws.Range("M2").Value = "X"
$ws.Range("N2").Value = "Y"
$ws.Range("O2").Value = "X"
$ws.Range("P2").Value = "Y"

# NGC4258 data (columns M:N), rows 3-7
$ws.Range("M3").Value = -20.928427534010702
$ws.Range("N3").Value = 0.99827686625497203
$ws.Range("M4").Value = -14.326310452583501
$ws.Range("N4").Value = 2.9783333588174199
$ws.Range("M5").Value = -12.938960891248399
$ws.Range("N5").Value = 3.9829083134689802
$ws.Range("M6").Value = -11.5335210985227
$ws.Range("N6").Value = 4.9803238053782799
$ws.Range("M7").Value = -10.737072829374601
$ws.Range("N7").Value = 7.0453599279374099

# NGC4631 data (columns O:P), rows 3-11
$ws.Range("O3").Value = -20.2448234448165
$ws.Range("P3").Value = 0.99850609175770999
$ws.Range("O4").Value = -18.8644138699485
$ws.Range("P4").Value = 1.9927097222925501
$ws.Range("O5").Value = -16.7281163912841
$ws.Range("P5").Value = 3.0000166045424201
$ws.Range("O6").Value = -13.724150960127799
$ws.Range("P6").Value = 4.0033132083253697
$ws.Range("O7").Value = -13.602778307915599
$ws.Range("P7").Value = 5.00368019170672
$ws.Range("O8").Value = -12.8719360441691
$ws.Range("P8").Value = 6.0082084910268501
$ws.Range("O9").Value = -10.635486849689499
$ws.Range("P9").Value = 7.00784056674732
$ws.Range("O10").Value = -10.4123277292942
$ws.Range("P10").Value = 9.0220426825113993
$ws.Range("O11").Value = -9.5987916712451309
$ws.Range("P11").Value = 10.0214932825539

# bestFit-ish column widths for the newly filled columns (M:P)
$ws.Columns.Item(13).ColumnWidth = 12
$ws.Columns.Item(14).ColumnWidth = 11.33
$ws.Columns.Item(15).ColumnWidth = 12
$ws.Columns.Item(16).ColumnWidth = 11.33

# Defined name that Excel creates for the imported "wpd_datasets (3)" CSV range,
# scoped to this worksheet.
$ws.Names.Add("wpd_datasets__3", "=Sheet1!`$M`$1:`$P`$11")

# Selection moved to the newly-imported data, as left by the import.
$ws.Range("M2").Select()
